$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.507.89'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '2.448.07'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'509.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").Value = "'132.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = "'0.557"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("D9").Value = '2.450.08'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = "'4.59"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -7.68%  '
$ws.Range("D14").Value = '2.876.82'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").Value = '57.517.24'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = "'21.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '2.423.30'
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").Value = "'4.10"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = "'313.98"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("D22").Value = "'6.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.57%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("D25").Value = "'65.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.27%  '
$ws.Range("D26").Value = "'0.994"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").Value = '2.522.49'
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("E28").Value = '  -5.46%  '
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").Value = "'7.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.27%  '
$ws.Range("D31").Value = "'173.84"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("D34").Value = "'6.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = "'0.996"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = "'17.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("D39").Value = "'1.24"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.81%  '
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("D42").Value = "'0.814"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = "'134.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +10.39%  '
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E46").Value = '  +4.09%  '
$ws.Range("D47").Value = "'255.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").Value = "'0.0918"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("D50").Value = "'0.0492"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("E51").Value = '  +0.87%  '
